# Update net value for Explorer — populate the new weekly row (row 35,
# 2017-04-28) that was appended to the "fund record" worksheet, and move
# the saved selection down to where the new data was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 35 already carries placeholder cells/styles inherited from the
# blank template rows below it (B35,C35,D35,I35,M35 have styles but no
# values). Most of the row just needs its value filled in using the
# formatting that is already in place; a handful of cells are brand new
# (E:H, A, P:U) and need both a value and the matching number format
# copied from the row above (row 34), which already has the exact same
# layout for this weekly table.
# ---------------------------------------------------------------------

# Date for the new week
$ws.Range("A35").Value = 20170428

# Total equity = SUM(C35:H35), same formula pattern as B29:B34
$ws.Range("B35").Formula = "=SUM(C35:H35)"

# Existing-style numeric cells (style already present on the row) -- just
# set the values, keep the inherited "#,##0.00" / "#,##0.0000_ " formats.
$ws.Range("C35").Value = 1122370
$ws.Range("D35").Value = 1164701
$ws.Range("M35").Value = 0.6822

# New cells that need the "0.00_ " format used by E34:H34.
$ws.Range("E35").Value = 844724
$ws.Range("F35").Value = 402437
$ws.Range("G35").Value = 599975
$ws.Range("H35").Value = 971971
"E35", "F35", "G35", "H35" | ForEach-Object {
    $ws.Range($_).NumberFormat = "0.00_ "
}

# N35 / O35 already hold placeholder styles for a different column in the
# blank rows below (N uses "#,##0.0_ ", O uses the same) -- re-point them
# at the formats actually used in row 34 for these columns.
$ws.Range("N35").Value = 2.1203
$ws.Range("N35").NumberFormat = "#,##0.0000_ "

$ws.Range("O35").Value = 424060
$ws.Range("O35").NumberFormat = "#,##0.00_ "

# Remaining brand-new cells, formatted like row 34.
$ws.Range("P35").Value = 1888524
$ws.Range("P35").NumberFormat = "#,##0.00"

$ws.Range("Q35").Value = 666684.18
$ws.Range("R35").Value = 2.8327

$ws.Range("S35").Value = 2793594
$ws.Range("S35").NumberFormat = "#,##0.00"

$ws.Range("T35").Value = 1748728.62
$ws.Range("T35").NumberFormat = "#,##0.00"

$ws.Range("U35").Value = 1.5975
$ws.Range("U35").NumberFormat = '0.0000_);[Red]\(0.0000\)'

# I35 was an empty placeholder cell for a column that isn't part of this
# week's data any more -- drop it entirely instead of leaving a styled
# blank cell behind.
$ws.Range("I35").Clear()

# Move the saved cursor/selection to where the new row was typed, and
# scroll the sheet up so row 16 is back at the top of the view.
$ws.Range("S35").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
